$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update balance on row 6 (F6: 1008 -> 1000)
$ws.Range("F6").Value = 1000

# Append 8 more rows (93-100) that mirror the existing "moses/bro" rows
# (90-92) by copying row 92 downward, preserving cell types/styles
# (notably keeping the numeric-looking ID "1234" stored as text).
$ws.Range("A92:F92").Copy()
for ($r = 93; $r -le 100; $r++) {
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial()
}
